$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")
$elements = $wb.Worksheets.Item("Elements")

$newUrl = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/employee-retirement-date"

$meta.Range("B2").Value = $newUrl
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

# Extension.url row also carries the same URL as its Fixed Value (shared string)
$elements.Range("Q5").Value = $newUrl
